# Generate Report for Handoff
# Re-order the localization status rows so the newly-handed-off file
# (c6aa9706...) appears first and the previously in-flight file
# (3f9fb24c...) moves to row 3 with a fresh "Ready for handoff" status
# and an updated handoff timestamp.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Overview" ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = "c6aa9706-a694-448a-8730-9f92d51da86d.md"
$ws.Range("A3").Value = "3f9fb24c-623b-4e25-9556-cc1b30c165f3.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

# ---- Sheet "zh-cn" ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = "c6aa9706-a694-448a-8730-9f92d51da86d.md"
$ws.Range("C2").Value = "c6aa9706-a694-448a-8730-9f92d51da86d.b0ff16d37518f09d7a603c08af57cac343352032.zh-cn.xlf"
$ws.Range("A3").Value = "3f9fb24c-623b-4e25-9556-cc1b30c165f3.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "3f9fb24c-623b-4e25-9556-cc1b30c165f3.34fdea629a3d59e361c8c887451dceb8dd979798.zh-cn.xlf"
$ws.Range("D3").Value = "2016-03-10 02:47:25"

# ---- Sheet "de-de" ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = "c6aa9706-a694-448a-8730-9f92d51da86d.md"
$ws.Range("C2").Value = "c6aa9706-a694-448a-8730-9f92d51da86d.b0ff16d37518f09d7a603c08af57cac343352032.de-de.xlf"
$ws.Range("A3").Value = "3f9fb24c-623b-4e25-9556-cc1b30c165f3.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "3f9fb24c-623b-4e25-9556-cc1b30c165f3.34fdea629a3d59e361c8c887451dceb8dd979798.de-de.xlf"
$ws.Range("D3").Value = "2016-03-10 02:47:28"

$wb.Save()
